$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force Text format first so decimal-looking values
# (e.g. "604.25", "1.00") are kept as literal text instead of being
# auto-coerced into numbers by Excel type inference, matching how the rest
# of the column already stores prices as text.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D11", "D12", "D13", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.180.45"
$ws.Range("D3").Value = "3.555.31"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "604.25"
$ws.Range("D6").Value = "143.48"
$ws.Range("D7").Value = "3.554.91"
$ws.Range("D11").Value = "7.82"
$ws.Range("D12").Value = "0.411"
$ws.Range("D13").Value = "4.151.70"
$ws.Range("D15").Value = "30.00"
$ws.Range("D16").Value = "3.542.36"
$ws.Range("D17").Value = "66.225.81"
$ws.Range("D19").Value = "11.33"
$ws.Range("D20").Value = "6.18"
$ws.Range("D21").Value = "14.65"
$ws.Range("D22").Value = "428.84"
$ws.Range("D23").Value = "0.608"
$ws.Range("D24").Value = "79.78"
$ws.Range("D25").Value = "3.696.04"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("D28").Value = "2.49"
$ws.Range("D29").Value = "9.08"
$ws.Range("D30").Value = "7.84"
$ws.Range("D32").Value = "3.548.96"
$ws.Range("D33").Value = "25.35"
$ws.Range("D34").Value = "1.44"
$ws.Range("D37").Value = "7.82"
$ws.Range("D39").Value = "5.54"
$ws.Range("D40").Value = "174.27"
$ws.Range("D41").Value = "0.0847"
$ws.Range("D42").Value = "5.18"
$ws.Range("D43").Value = "0.885"
$ws.Range("D44").Value = "1.91"
$ws.Range("D45").Value = "45.95"
$ws.Range("D48").Value = "24.78"
$ws.Range("D50").Value = "7.09"
$ws.Range("D51").Value = "22.83"

# --- Volume(1h) column (E): percentage strings, always text already.
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").Value = "  +7.57%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("E35").Value = "  -8.97%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("E51").Value = "  +1.61%  "
